$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated pension-recipient figures for 2015-2021 (row 4)
$ws.Range("E4").Value = 12138
$ws.Range("F4").Value = 12326
$ws.Range("G4").Value = 12492
$ws.Range("H4").Value = 12672
$ws.Range("I4").Value = 12953
$ws.Range("J4").Value = 13314
$ws.Range("K4").Value = 13522

# Reflect the author's active selection at save time
$ws.Range("E4:K4").Select()
